# The target edit only touches word/styles.xml's <w:docDefaults> block
# (<w:rPrDefault>/<w:rPr> and <w:pPrDefault>/<w:pPr>), trimming a set of
# explicit, schema-default-valued properties (b=0, i=0, smallCaps=0,
# strike=0, color=000000, u=none, shd=clear/auto, vertAlign=baseline,
# keepNext=0, keepLines=0, widowControl=1, the empty pBdr, spacing
# before/after=0, ind=0/0/0, contextualSpacing=0, jc=left) down to a
# minimal remainder (rFonts/sz/szCs/lang for run defaults, and a bare
# line-spacing <w:spacing w:line="276" w:lineRule="auto"/> for paragraph
# defaults). No <w:style> element (Normal, headings, ...) changes at all.
#
# docDefaults is a document-creation-time fallback that Word's object
# model does not expose for editing: there is no Styles item, property,
# or method (in real Word's VBA/COM surface, and equally none here) that
# reads or writes <w:docDefaults>/rPrDefault/pPrDefault. Style objects
# (Styles.Item(...).Font / .ParagraphFormat) only ever read/write the
# explicit <w:rPr>/<w:pPr> carried on that individual <w:style> element
# (e.g. "Normal"), never the package-level docDefaults fallback, and
# Document/Range WordOpenXML and InsertXML only round-trip body content
# (word/document.xml) - confirmed read-only / content-only in this host.
#
# Since docDefaults can't be reached from the Word COM object model,
# there is no automation call that reproduces this hunk without forging
# an unrequested change elsewhere (e.g. stamping redundant formatting
# onto the Normal style, which the diff does NOT do). To avoid
# introducing edits that aren't actually part of the described change,
# this script intentionally performs no Word object-model mutation.

$d = $word.ActiveDocument
# Touch the document object without mutating any content/formatting.
$null = $d.Name
